$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.675.49'
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').Value = '3.494.25'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.12'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.73'
$ws.Range('E6').Value = '  +4.39%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.496.69'
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.593'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  +6.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.10'
$ws.Range('E11').Value = '  -2.20%  '
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').Value = '4.103.82'
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.23'
$ws.Range('E14').Value = '  +11.78%  '
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = '67.645.78'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '3.498.22'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.28'
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.37'
$ws.Range('E20').Value = '  +2.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '393.02'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').Value = '  -1.04%  '
$ws.Range('E23').Value = '  +0.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000121'
$ws.Range('E27').Value = '  +1.37%  '
$ws.Range('E28').Value = '  +2.30%  '
$ws.Range('E29').Value = '  -2.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.21'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.43'
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.06'
$ws.Range('E33').Value = '  +0.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.65'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.40'
$ws.Range('E35').Value = '  +0.68%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '163.98'
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.876'
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.89'
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.74'
$ws.Range('E41').Value = '  +7.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.88'
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.65'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('D44').Value = '2.854.06'
$ws.Range('E44').Value = '  +1.72%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '26.19'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0726'
$ws.Range('E46').Value = '  -2.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.67'
$ws.Range('E47').Value = '  -2.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '42.10'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0300'
$ws.Range('E49').Value = '  -0.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '337.82'
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('E51').Value = '  -1.37%  '
